$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 for columns A, B, E, F, G, I, M
# (all other columns share identical values between the two rows)

$cols = @("A", "B", "E", "F", "G", "I", "M")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}
